$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose updated "Price" text would otherwise be auto-parsed
# as a number by Excel (because it looks like a plain decimal). Column D in this
# sheet is always text (e.g. "27.140.67", "336.58", ...), so force a text number
# format before assigning the value, matching the original inline-string cells.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.140.67'
$ws.Range("E2").Value = '  -0.87%  '

# Row 3
$ws.Range("D3").Value = '1.781.93'
$ws.Range("E3").Value = '  -1.69%  '

# Row 4
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '336.58'
$ws.Range("E5").Value = '  -2.03%  '

# Row 6
$ws.Range("E6").Value = '  +0.04%  '

# Row 7
$ws.Range("D7").Value = '0.3828'
$ws.Range("E7").Value = '  +0.45%  '

# Row 8
$ws.Range("D8").Value = '0.3421'
$ws.Range("E8").Value = '  -1.95%  '

# Row 9
$ws.Range("D9").Value = '47.93'
$ws.Range("E9").Value = '  -2.01%  '

# Row 10
$ws.Range("D10").Value = '1.187'
$ws.Range("E10").Value = '  -3.45%  '

# Row 11
$ws.Range("D11").Value = '0.07456'
$ws.Range("E11").Value = '  -3.45%  '

# Row 12
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.19%  '

# Row 13
$ws.Range("D13").Value = '21.65'
$ws.Range("E13").Value = '  -1.76%  '

# Row 14
$ws.Range("D14").Value = '6.436'
$ws.Range("E14").Value = '  -2.37%  '

# Row 15
$ws.Range("D15").Value = '1.782.00'
$ws.Range("E15").Value = '  -1.76%  '

# Row 16
$ws.Range("D16").Value = '7.103'
$ws.Range("E16").Value = '  -1.37%  '

# Row 17
$ws.Range("D17").Value = '0.00001093'
$ws.Range("E17").Value = '  -2.08%  '

# Row 18
$ws.Range("D18").Value = '0.06665'
$ws.Range("E18").Value = '  -0.80%  '

# Row 19
$ws.Range("D19").Value = '83.47'
$ws.Range("E19").Value = '  -2.93%  '

# Row 20
$ws.Range("E20").Value = '  +0.01%  '

# Row 21
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '17.41'
$ws.Range("E21").Value = '  -0.85%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.507'
$ws.Range("E22").Value = '  -0.64%  '

# Row 23
$ws.Range("D23").Value = '27.133.19'
$ws.Range("E23").Value = '  -0.93%  '

# Row 24
$ws.Range("D24").Value = '12.25'
$ws.Range("E24").Value = '  -7.24%  '

# Row 25
$ws.Range("E25").Value = '  -3.84%  '

# Row 26
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.502'
$ws.Range("E26").Value = '  -5.67%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '21.17'
$ws.Range("E27").Value = '  -3.73%  '

# Row 28
$ws.Range("D28").Value = '1.441'
$ws.Range("E28").Value = '  -1.41%  '

# Row 29
$ws.Range("D29").Value = '154.73'
$ws.Range("E29").Value = '  +0.56%  '

# Row 30
$ws.Range("D30").Value = '1.985.66'
$ws.Range("E30").Value = '  -1.65%  '

# Row 31
$ws.Range("D31").Value = '134.19'
$ws.Range("E31").Value = '  -1.06%  '

# Row 32
$ws.Range("D32").Value = '3.975'
$ws.Range("E32").Value = '  -1.45%  '

# Row 33
$ws.Range("D33").Value = '6.006'
$ws.Range("E33").Value = '  -4.58%  '

# Row 34
$ws.Range("D34").Value = '0.08670'
$ws.Range("E34").Value = '  -1.19%  '

# Row 35
$ws.Range("D35").Value = '13.03'
$ws.Range("E35").Value = '  -5.90%  '

# Row 36
$ws.Range("D36").Value = '1.623'
$ws.Range("E36").Value = '  -3.87%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '5.394'
$ws.Range("E37").Value = '  -3.65%  '

# Row 38
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '0.6820'
$ws.Range("E38").Value = '  -1.66%  '

# Row 39
$ws.Range("D39").Value = '0.06329'
$ws.Range("E39").Value = '  -2.00%  '

# Row 40
$ws.Range("D40").Value = '0.02331'
$ws.Range("E40").Value = '  -2.65%  '

# Row 41
$ws.Range("D41").Value = '0.2178'
$ws.Range("E41").Value = '  -3.89%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '8.437'
$ws.Range("E42").Value = '  -5.13%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.236'
$ws.Range("E43").Value = '  -4.79%  '

# Row 44
$ws.Range("D44").Value = '14.19'
$ws.Range("E44").Value = '  -3.56%  '

# Row 45
$ws.Range("E45").Value = '  -0.04%  '

# Row 46
$ws.Range("D46").Value = '0.6410'
$ws.Range("E46").Value = '  -1.27%  '

# Row 47
$ws.Range("D47").Value = '3.856'
$ws.Range("E47").Value = '  -4.04%  '

# Row 48
$ws.Range("D48").Value = '2.176'
$ws.Range("E48").Value = '  +0.40%  '

# Row 49
$ws.Range("D49").Value = '131.24'
$ws.Range("E49").Value = '  -0.99%  '

# Row 50
$ws.Range("D50").Value = '0.07099'
$ws.Range("E50").Value = '  -3.00%  '

# Row 51
$ws.Range("D51").Value = '78.61'
$ws.Range("E51").Value = '  -2.00%  '

# Restore the default (Normal) style on the cells we temporarily reformatted
# above, now that their text values are committed, so this script does not
# leave behind any extra per-cell formatting that was not in the original file.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
